# Scheduled market-data refresh for the Leve profit-tracking workbook.
# Updates currentAveragePrice(NQ/HQ) and recalculated Leve profit figures
# on each job sheet (Table_<JOB>) to reflect the latest market board data.

$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 323.25
$ws.Range("I9").Value = 449
$ws.Range("K9").Value = 449
$ws.Range("M9").Value = -280
$ws.Range("H21").Value = 24999.055
$ws.Range("I21").Value = 20000
$ws.Range("J21").Value = 25293.117
$ws.Range("K21").Value = 20000
$ws.Range("L21").Value = 25293.117
$ws.Range("M21").Value = -19532
$ws.Range("N21").Value = -26229.117
$ws.Range("H23").Value = 24999.055
$ws.Range("I23").Value = 20000
$ws.Range("J23").Value = 25293.117
$ws.Range("K23").Value = 20000
$ws.Range("L23").Value = 25293.117
$ws.Range("M23").Value = -19766
$ws.Range("N23").Value = -25761.117
$ws.Range("H38").Value = 3067.1428
$ws.Range("I38").Value = 3067.1428
$ws.Range("K38").Value = 9201.428400000001
$ws.Range("M38").Value = -8829.428400000001
$ws.Range("H41").Value = 1176.6
$ws.Range("I41").Value = 1843.6
$ws.Range("K41").Value = 1843.6
$ws.Range("M41").Value = -1403.6
$ws.Range("H43").Value = 2347.5
$ws.Range("I43").Value = 2517
$ws.Range("K43").Value = 2517
$ws.Range("M43").Value = -2448
$ws.Range("H76").Value = 4424.143
$ws.Range("I76").Value = 2984.5
$ws.Range("K76").Value = 2984.5
$ws.Range("M76").Value = -2669.5
$ws.Range("H79").Value = 4424.143
$ws.Range("I79").Value = 2984.5
$ws.Range("K79").Value = 2984.5
$ws.Range("M79").Value = -1892.5
$ws.Range("H99").Value = 680.2
$ws.Range("I99").Value = 311.33334
$ws.Range("J99").Value = 4000
$ws.Range("K99").Value = 934.0000200000001
$ws.Range("L99").Value = 12000
$ws.Range("M99").Value = 563.9999799999999
$ws.Range("N99").Value = -14996
$ws.Range("H132").Value = 56123.26
$ws.Range("I132").Value = 69982.87
$ws.Range("K132").Value = 209948.61
$ws.Range("M132").Value = -207418.61
$ws.Range("H137").Value = 3636.3333
$ws.Range("J137").Value = 4823.625
$ws.Range("L137").Value = 14470.875
$ws.Range("N137").Value = -19570.875

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 8520.789000000001
$ws.Range("I2").Value = 9632.714
$ws.Range("J2").Value = 5407.4
$ws.Range("K2").Value = 9632.714
$ws.Range("L2").Value = 5407.4
$ws.Range("M2").Value = -9519.714
$ws.Range("N2").Value = -5633.4
$ws.Range("H61").Value = 3226.8572
$ws.Range("I61").Value = 3226.8572
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3226.8572
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -3014.8572
$ws.Range("H116").Value = 8520.789000000001
$ws.Range("I116").Value = 9632.714
$ws.Range("J116").Value = 5407.4
$ws.Range("K116").Value = 9632.714
$ws.Range("L116").Value = 5407.4
$ws.Range("M116").Value = -7338.714
$ws.Range("N116").Value = -9995.4
$ws.Range("H136").Value = 3226.8572
$ws.Range("I136").Value = 3226.8572
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 9680.571599999999
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -7130.571599999999
$ws.Range("N61").ClearContents()
$ws.Range("N136").ClearContents()

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 8520.789000000001
$ws.Range("I3").Value = 9632.714
$ws.Range("J3").Value = 5407.4
$ws.Range("K3").Value = 9632.714
$ws.Range("L3").Value = 5407.4
$ws.Range("M3").Value = -9518.714
$ws.Range("N3").Value = -5635.4
$ws.Range("H5").Value = 124.2
$ws.Range("I5").Value = 119.57143
$ws.Range("K5").Value = 119.57143
$ws.Range("M5").Value = -6.571430000000007
$ws.Range("H86").Value = 1799.1428
$ws.Range("I86").Value = 1758.8
$ws.Range("K86").Value = 1758.8
$ws.Range("M86").Value = -635.8
$ws.Range("H89").Value = 1799.1428
$ws.Range("I89").Value = 1758.8
$ws.Range("K89").Value = 8794
$ws.Range("M89").Value = -3178

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 25127
$ws.Range("J59").Value = 25127
$ws.Range("L59").Value = 25127
$ws.Range("N59").Value = -27417

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 86637.03
$ws.Range("I4").Value = 919.76
$ws.Range("K4").Value = 2759.28
$ws.Range("M4").Value = -2647.28
$ws.Range("H5").Value = 1899.5
$ws.Range("I5").Value = 2849.5
$ws.Range("J5").Value = 949.5
$ws.Range("K5").Value = 8548.5
$ws.Range("L5").Value = 2848.5
$ws.Range("M5").Value = -8436.5
$ws.Range("N5").Value = -3072.5
$ws.Range("H7").Value = 498.93332
$ws.Range("I7").Value = 123.55556
$ws.Range("K7").Value = 370.66668
$ws.Range("M7").Value = -258.66668
$ws.Range("H39").Value = 4562.857
$ws.Range("J39").Value = 4837.077
$ws.Range("L39").Value = 14511.231
$ws.Range("N39").Value = -15099.231
$ws.Range("H121").Value = 1401747.9
$ws.Range("I121").Value = 112467.78
$ws.Range("K121").Value = 337403.34
$ws.Range("M121").Value = -336093.34
$ws.Range("H135").Value = 1899.5
$ws.Range("I135").Value = 2849.5
$ws.Range("J135").Value = 949.5
$ws.Range("K135").Value = 25645.5
$ws.Range("L135").Value = 8545.5
$ws.Range("M135").Value = -23110.5
$ws.Range("N135").Value = -13615.5

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 19989.25
$ws.Range("J5").Value = 19989.25
$ws.Range("L5").Value = 19989.25
$ws.Range("N5").Value = -20213.25
$ws.Range("H33").Value = 21100
$ws.Range("J33").Value = 21100
$ws.Range("L33").Value = 21100
$ws.Range("N33").Value = -21604
$ws.Range("H44").Value = 27000
$ws.Range("I44").Value = 30000
$ws.Range("J44").Value = 25500
$ws.Range("K44").Value = 30000
$ws.Range("L44").Value = 25500
$ws.Range("M44").Value = -29404
$ws.Range("N44").Value = -26692
$ws.Range("H46").Value = 10041
$ws.Range("I46").Value = 10041
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 10041
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -9885
$ws.Range("H52").Value = 15699.667
$ws.Range("J52").Value = 15699.667
$ws.Range("L52").Value = 15699.667
$ws.Range("N52").Value = -16217.667
$ws.Range("H57").Value = 13309.833
$ws.Range("I57").Value = 4554
$ws.Range("K57").Value = 4554
$ws.Range("M57").Value = -3734
$ws.Range("H122").Value = 4170.0835
$ws.Range("I122").Value = 6392
$ws.Range("J122").Value = 2583
$ws.Range("K122").Value = 19176
$ws.Range("L122").Value = 7749
$ws.Range("M122").Value = -16726
$ws.Range("N122").Value = -12649
$ws.Range("H132").Value = 2899.5
$ws.Range("J132").Value = 2999.5
$ws.Range("L132").Value = 8998.5
$ws.Range("N132").Value = -14058.5
$ws.Range("N46").ClearContents()

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 10137.272
$ws.Range("J43").Value = 10151.2
$ws.Range("L43").Value = 10151.2
$ws.Range("N43").Value = -10537.2
$ws.Range("H61").Value = 2901.0833
$ws.Range("I61").Value = 1881.4
$ws.Range("K61").Value = 1881.4
$ws.Range("M61").Value = -1679.4
$ws.Range("H113").Value = 2901.0833
$ws.Range("I113").Value = 1881.4
$ws.Range("K113").Value = 1881.4
$ws.Range("M113").Value = 288.5999999999999

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 18022.375
$ws.Range("J45").Value = 23424.25
$ws.Range("L45").Value = 23424.25
$ws.Range("N45").Value = -24406.25
$ws.Range("H81").Value = 7091.933
$ws.Range("I81").Value = 7091.933
$ws.Range("K81").Value = 14183.866
$ws.Range("M81").Value = -13122.866
$ws.Range("H84").Value = 7091.933
$ws.Range("I84").Value = 7091.933
$ws.Range("K84").Value = 70919.33
$ws.Range("M84").Value = -65615.33
$ws.Range("H122").Value = 3637.65
$ws.Range("I122").Value = 4619.6
$ws.Range("K122").Value = 13858.8
$ws.Range("M122").Value = -11408.8
